# Add a new "Request" type: estimation across multiple Synopses.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Request")

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "ESTIMATE"
$ws.Range("C9").Value = "request an estimation among multiple Synopses"

# Row 8 (B column) already carries the bold/centered "OperationType" style;
# copy that formatting onto the newly-filled B9 cell (matches what Excel
# does automatically when you type into a cell that inherits row style).
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Resize the structured table (Table16) so it covers the new row too.
$table = $ws.ListObjects.Item("Table16")
$table.Resize($ws.Range("A1:C9"))

# Make "Request" the active sheet/tab, mirroring the saved selection state.
$ws.Activate()
$ws.Range("A10").Select()
